# Generate Report for Handoff
#
# Swap the two data rows on each sheet so that the file that is
# "Handed back: in sync with en-US" (388cb62e...) sits in row 2 and the
# file that has just become "Ready for handoff" (131a7ed5...) sits in
# row 3, and refresh its handoff status/date fields.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "131a7ed5-056e-4060-bc28-41101af30063.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-22 12:52:56"

$ws.Range("A3").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-22 12:54:37"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
    }
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-22 12:52:52"
$ws.Range("F2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
$ws.Range("G2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-22 12:53:49"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "131a7ed5-056e-4060-bc28-41101af30063.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-22 12:54:33"
$ws.Range("F3").Value = "131a7ed5-056e-4060-bc28-41101af30063.md"
$ws.Range("G3").Value = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-22 12:53:49"
$ws.Range("J3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
    }
    if ($addr -eq '$D$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.zh-cn.xlf"
    }
    if ($addr -eq '$F$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
    }
    if ($addr -eq '$G$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.zh-cn.xlf"
    }
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.md"
    }
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf"
    }
    if ($addr -eq '$F$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.md"
    }
    if ($addr -eq '$G$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.de-de.xlf"
$ws.Range("E2").Value = "2016-03-22 12:52:56"
$ws.Range("F2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
$ws.Range("G2").Value = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.de-de.xlf"
$ws.Range("H2").Value = "2016-03-22 12:53:58"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "131a7ed5-056e-4060-bc28-41101af30063.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf"
$ws.Range("E3").Value = "2016-03-22 12:54:37"
$ws.Range("F3").Value = "131a7ed5-056e-4060-bc28-41101af30063.md"
$ws.Range("G3").Value = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf"
$ws.Range("H3").Value = "2016-03-22 12:53:58"
$ws.Range("J3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
    }
    if ($addr -eq '$D$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.de-de.xlf"
    }
    if ($addr -eq '$F$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.md"
    }
    if ($addr -eq '$G$2') {
        $hl.TextToDisplay = "388cb62e-4b38-4fce-8e47-1c7b6e786d71.e14923f95111d2ddee50c6ec3306a0c813bf63af.de-de.xlf"
    }
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.md"
    }
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf"
    }
    if ($addr -eq '$F$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.md"
    }
    if ($addr -eq '$G$3') {
        $hl.TextToDisplay = "131a7ed5-056e-4060-bc28-41101af30063.943fe33fac80d81fc32621ff878cd0d36d500c4f.de-de.xlf"
    }
}
